# "create the 4, 8 and 12 version"
# The worksheet tracked a single child ("Luis") practicing 4 failed words.
# It is edited to track a different child ("Paco") and a bigger 12-word
# failure list (rows 2..13), updating the two summary mini-tables (G:H and
# K:L) and the header/score cells (B2:D2, J2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name + time/score cells -------------------------------------------------
$ws.Range("B2").Value2 = "Paco"
$ws.Range("C2").Value2 = 5
$ws.Range("D2").Value2 = 60
$ws.Range("J2").Value2 = "Paco:"

# --- Word list (col G) + fails (col H), and the mirrored col K/L table ------
# col L stores "0"/"1" as TEXT (shared string), matching the source file, so
# the target range is pre-formatted as Text before the values are poured in,
# then restored to the default style (no visible format change, just avoids
# Excel's normal "looks like a number" auto-conversion).
$lRange = $ws.Range("L2:L13")
$lRange.NumberFormat = "@"

# "Fail" here is the text ("0"/"1") that lands in column L (col H, the
# numeric fails-per-word tally, is reset to 0 for every word, per the diff).
$rows = @(
  @{ Row = 2;  Word = "falda";    Fail = 0 },
  @{ Row = 3;  Word = "aporbar";  Fail = 1 },
  @{ Row = 4;  Word = "abarzar";  Fail = 1 },
  @{ Row = 5;  Word = "plamera";  Fail = 1 },
  @{ Row = 6;  Word = "furta";    Fail = 1 },
  @{ Row = 7;  Word = "tornillo"; Fail = 0 },
  @{ Row = 8;  Word = "galdiolo"; Fail = 1 },
  @{ Row = 9;  Word = "palmera";  Fail = 0 },
  @{ Row = 10; Word = "parque";   Fail = 0 },
  @{ Row = 11; Word = "tractor";  Fail = 0 },
  @{ Row = 12; Word = "flor";     Fail = 0 },
  @{ Row = 13; Word = "adrono";   Fail = 1 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $word = $r.Word
  $fail = $r.Fail

  $ws.Range("G$row").Value2 = $word
  $ws.Range("H$row").Value2 = 0
  $ws.Range("K$row").Value2 = $word
  $ws.Range("L$row").Value2 = [string]$fail
}

$lRange.Style = "Normal"

# --- Chart 3's title mirrors the child's name --------------------------------
$chartObj = $ws.ChartObjects(3)
$chartObj.Chart.ChartTitle.Text = "Paco"

Write-Host "Updated Sheet1 B2:L13 for Paco (words through row 13) and chart title."
